$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 8211.111000000001
$ws.Range("J17").Value = 8612.5
$ws.Range("L17").Value = 25837.5
$ws.Range("N17").Value = -26173.5
$ws.Range("H28").Value = 1338.32
$ws.Range("I28").Value = 977.35
$ws.Range("J28").Value = 2782.2
$ws.Range("K28").Value = 977.35
$ws.Range("L28").Value = 2782.2
$ws.Range("M28").Value = -492.35
$ws.Range("N28").Value = -3752.2
$ws.Range("H40").Value = 7776.5454
$ws.Range("I40").Value = 6466.6665
$ws.Range("K40").Value = 6466.6665
$ws.Range("M40").Value = -6291.6665
$ws.Range("H51").Value = 7647.129
$ws.Range("I51").Value = 7538.375
$ws.Range("K51").Value = 7538.375
$ws.Range("M51").Value = -7054.375
$ws.Range("H113").Value = 8230.951999999999
$ws.Range("I113").Value = 5977.3335
$ws.Range("J113").Value = 9132.4
$ws.Range("K113").Value = 5977.3335
$ws.Range("L113").Value = 9132.4
$ws.Range("M113").Value = -2723.3335
$ws.Range("N113").Value = -15640.4
$ws.Range("H116").Value = 18188
$ws.Range("J116").Value = 18502
$ws.Range("L116").Value = 18502
$ws.Range("N116").Value = -25386
$ws.Range("H132").Value = 1971.0303
$ws.Range("I132").Value = 1905.4839
$ws.Range("K132").Value = 5716.4517
$ws.Range("M132").Value = -3186.4517
$ws.Range("H137").Value = 2859.5068
$ws.Range("I137").Value = 1634.0385
$ws.Range("K137").Value = 4902.1155
$ws.Range("M137").Value = -2352.1155
$ws.Range("H138").Value = 2554.342
$ws.Range("I138").Value = 1212.1724
$ws.Range("J138").Value = 3382.4893
$ws.Range("K138").Value = 3636.5172
$ws.Range("L138").Value = 10147.4679
$ws.Range("M138").Value = 1503.4828
$ws.Range("N138").Value = -20427.4679

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 20838440
$ws.Range("I74").Value = 30306834
$ws.Range("K74").Value = 30306834
$ws.Range("M74").Value = -30305960
$ws.Range("H77").Value = 20838440
$ws.Range("I77").Value = 30306834
$ws.Range("K77").Value = 151534170
$ws.Range("M77").Value = -151529802
$ws.Range("H97").Value = 813.4400000000001
$ws.Range("I97").Value = 748
$ws.Range("K97").Value = 748
$ws.Range("M97").Value = -252
$ws.Range("H108").Value = 90000
$ws.Range("J108").Value = 90000
$ws.Range("L108").Value = 90000
$ws.Range("N108").Value = -97680
$ws.Range("H132").Value = 7747.4614
$ws.Range("I132").Value = 4671.2856
$ws.Range("K132").Value = 14013.8568
$ws.Range("M132").Value = -11483.8568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1258.5
$ws.Range("I134").Value = 1120.9445
$ws.Range("K134").Value = 3362.8335
$ws.Range("M134").Value = -827.8335000000002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 985
$ws.Range("I25").Value = 947.5
$ws.Range("K25").Value = 947.5
$ws.Range("M25").Value = -773.5
$ws.Range("H31").Value = 27606.4
$ws.Range("I31").Value = 2855.7812
$ws.Range("J31").Value = 88531
$ws.Range("K31").Value = 2855.7812
$ws.Range("L31").Value = 88531
$ws.Range("M31").Value = -2560.7812
$ws.Range("N31").Value = -89121
$ws.Range("H34").Value = 27606.4
$ws.Range("I34").Value = 2855.7812
$ws.Range("J34").Value = 88531
$ws.Range("K34").Value = 2855.7812
$ws.Range("L34").Value = 88531
$ws.Range("M34").Value = -2653.7812
$ws.Range("N34").Value = -88935
$ws.Range("H58").Value = 2832
$ws.Range("I58").Value = 1519.5264
$ws.Range("J58").Value = 6988.1665
$ws.Range("K58").Value = 1519.5264
$ws.Range("L58").Value = 6988.1665
$ws.Range("M58").Value = -1316.5264
$ws.Range("N58").Value = -7394.1665
$ws.Range("H132").Value = 2554.6365
$ws.Range("I132").Value = 1892.0333
$ws.Range("K132").Value = 5676.0999
$ws.Range("M132").Value = -3146.0999
$ws.Range("H136").Value = 2832
$ws.Range("I136").Value = 1519.5264
$ws.Range("J136").Value = 6988.1665
$ws.Range("K136").Value = 4558.5792
$ws.Range("L136").Value = 20964.4995
$ws.Range("M136").Value = -2008.5792
$ws.Range("N136").Value = -26064.4995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2567935
$ws.Range("J5").Value = 3337430.5
$ws.Range("L5").Value = 10012291.5
$ws.Range("N5").Value = -10012515.5
$ws.Range("H33").Value = 555618.2
$ws.Range("I33").Value = 714343.4
$ws.Range("J33").Value = 80
$ws.Range("K33").Value = 4286060.4
$ws.Range("L33").Value = 480
$ws.Range("M33").Value = -4285777.4
$ws.Range("N33").Value = -1046
$ws.Range("H46").Value = 2536.0527
$ws.Range("I46").Value = 29
$ws.Range("J46").Value = 3431.4285
$ws.Range("K46").Value = 87
$ws.Range("L46").Value = 10294.2855
$ws.Range("M46").Value = 4
$ws.Range("N46").Value = -10476.2855
$ws.Range("H117").Value = 3578.3333
$ws.Range("J117").Value = 3844
$ws.Range("L117").Value = 11532
$ws.Range("N117").Value = -18416
$ws.Range("H129").Value = 4904890.5
$ws.Range("I129").Value = 1246.7142
$ws.Range("J129").Value = 8337441
$ws.Range("K129").Value = 3740.1426
$ws.Range("L129").Value = 25012323
$ws.Range("M129").Value = 1259.8574
$ws.Range("N129").Value = -25022323
$ws.Range("H131").Value = 5917469.5
$ws.Range("I131").Value = 13890698
$ws.Range("K131").Value = 41672094
$ws.Range("M131").Value = -41667054
$ws.Range("H135").Value = 2567935
$ws.Range("J135").Value = 3337430.5
$ws.Range("L135").Value = 30036874.5
$ws.Range("N135").Value = -30041944.5
$ws.Range("H139").Value = 3753.6553
$ws.Range("I139").Value = 2529.9167
$ws.Range("K139").Value = 7589.750100000001
$ws.Range("M139").Value = -2449.750100000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4707.5127
$ws.Range("I122").Value = 3559.3447
$ws.Range("K122").Value = 10678.0341
$ws.Range("M122").Value = -8228.034100000001
$ws.Range("H132").Value = 5063.6665
$ws.Range("I132").Value = 4716.324
$ws.Range("K132").Value = 14148.972
$ws.Range("M132").Value = -11618.972
$ws.Range("H133").Value = 65977.60000000001
$ws.Range("J133").Value = 65977.60000000001
$ws.Range("L133").Value = 65977.60000000001
$ws.Range("N133").Value = -76097.60000000001
$ws.Range("H134").Value = 58485.2
$ws.Range("J134").Value = 58485.2
$ws.Range("L134").Value = 175455.6
$ws.Range("N134").Value = -180525.6
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
$ws.Range("H136").Value = 57587.848
$ws.Range("J136").Value = 57587.848
$ws.Range("L136").Value = 172763.544
$ws.Range("N136").Value = -177863.544

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5234.484
$ws.Range("I7").Value = 4494.729
$ws.Range("J7").Value = 7770.7856
$ws.Range("K7").Value = 4494.729
$ws.Range("L7").Value = 7770.7856
$ws.Range("M7").Value = -4382.729
$ws.Range("N7").Value = -7994.7856
$ws.Range("H63").Value = 39995
$ws.Range("I63").Value = 39995
$ws.Range("K63").Value = 39995
$ws.Range("M63").Value = -39246
$ws.Range("H66").Value = 39995
$ws.Range("I66").Value = 39995
$ws.Range("K66").Value = 119985
$ws.Range("M66").Value = -116241
$ws.Range("H122").Value = 186699.95
$ws.Range("I122").Value = 252905.06
$ws.Range("K122").Value = 758715.1799999999
$ws.Range("M122").Value = -756265.1799999999
$ws.Range("H126").Value = 5234.484
$ws.Range("I126").Value = 4494.729
$ws.Range("J126").Value = 7770.7856
$ws.Range("K126").Value = 13484.187
$ws.Range("L126").Value = 23312.3568
$ws.Range("M126").Value = -11014.187
$ws.Range("N126").Value = -28252.3568
$ws.Range("H132").Value = 13771.357
$ws.Range("I132").Value = 14055
$ws.Range("K132").Value = 42165
$ws.Range("M132").Value = -39635
$ws.Range("H136").Value = 6638.96
$ws.Range("J136").Value = 8069.1177
$ws.Range("L136").Value = 24207.3531
$ws.Range("N136").Value = -29307.3531
$ws.Range("H137").Value = 63630
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()
$ws.Range("H139").Value = 69990.836
$ws.Range("J139").Value = 69990.836
$ws.Range("L139").Value = 69990.836
$ws.Range("N139").Value = -80270.836

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1946.0312
$ws.Range("I126").Value = 1741.6923
$ws.Range("K126").Value = 5225.0769
$ws.Range("M126").Value = -2755.0769
$ws.Range("H132").Value = 2982.5117
$ws.Range("I132").Value = 2137.0967
$ws.Range("K132").Value = 6411.2901
$ws.Range("M132").Value = -3881.2901
$ws.Range("H136").Value = 2321.7114
$ws.Range("I136").Value = 1994.5883
$ws.Range("K136").Value = 5983.7649
$ws.Range("M136").Value = -3433.7649
$ws.Range("H137").Value = 69195
$ws.Range("J137").Value = 69195
$ws.Range("L137").Value = 69195
$ws.Range("N137").Value = -79395
